$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -------------------------------------------------------------
# Header row
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Total Reward"
$ws.Range("C1").Value = "Total Savings"
$ws.Range("D1").Value = "Total Tokens"
$ws.Range("E1").Value = "Total Time"

# steven's row (row 2) already has the username, add the metrics
$ws.Range("A2").Value = "steven"
$ws.Range("B2").Value = 10.0
$ws.Range("C2").Value = 5.0
$ws.Range("D2").Value = 1500.0
$ws.Range("E2").Value = 30.0

# connie's row (new row 3)
$ws.Range("A3").Value = "connie"
$ws.Range("B3").Value = 5.0
$ws.Range("C3").Value = 10.0
$ws.Range("D3").Value = 2000.0
$ws.Range("E3").Value = 30.0

# --- Formatting: match the style already used by column A (s="1") -------
$ws.Range("A1").Copy()
$ws.Range("A1:E3").PasteSpecial(-4122)

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 12.6
$ws.Columns.Item(4).ColumnWidth = 12.6
